$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing job title text
$ws.Range("A3").Value = "Java Developer"

# Add two new job titles
$ws.Range("A4").Value = "Data Engineer"
$ws.Range("A5").Value = "Fullstack Developer"

# Update the active selection to match the saved view state
$ws.Range("D9").Select() | Out-Null
